$d = $word.ActiveDocument

# Locate the byline paragraph "Edison Achalma" (style "Author") that sits
# right under the article title, then add a new "Author"-styled paragraph
# right after it with the author's institutional affiliation.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Edison Achalma" -and $p.Style.NameLocal -eq "Author") {
        $p.Range.InsertAfter("`r")
        $newPara = $p.Next()
        $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
        break
    }
}
